# Updated cryptos list on Fri Jun  2 06:29:50 UTC 2023 with GitHub Actions
#
# This script reproduces the refreshed "Price" (column D) and
# "Volume(1h)" (column E) figures scraped for each coin, plus the
# row-42/row-44 swap of FraxShare <-> Aptos.
#
# Because several of the new Price values look like plain decimal
# numbers (e.g. "0.9996", "1.000", "5.000"), assigning them straight
# to Range.Value would make Excel auto-convert the cell to a numeric
# type. The source workbook stores every Price/Volume cell as text,
# so for values that parse as a number we briefly mark the cell as
# Text (NumberFormat "@"), assign the value, then restore the cell's
# style to "Normal" so no stray number formatting/style index is left
# behind on the cell.

function Set-TextValue {
    param(
        $Cell,
        [string]$Value
    )
    if ($Value -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Value
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Value
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Price, Volume)  [$null means "leave unchanged"]
$rows = @{
    2  = @("27.138.53",   "  +1.09%  ")
    3  = @("1.891.45",    "  +1.97%  ")
    4  = @("0.9996",      "  -0.10%  ")
    5  = @("308.21",      "  +1.22%  ")
    6  = @("0.9998",      "  -0.04%  ")
    7  = @("0.5188",      "  +2.93%  ")
    8  = @("0.3723",      "  +2.00%  ")
    9  = @("0.07211",     "  +0.56%  ")
    10 = @("0.9051",      "  +1.75%  ")
    11 = @("21.07",       "  +2.01%  ")
    12 = @("0.07637",     "  +1.67%  ")
    13 = @("1.890.97",    "  +1.92%  ")
    14 = @("95.19",       "  +3.64%  ")
    15 = @("5.282",       "  +1.00%  ")
    16 = @("1.000",       "  -0.07%  ")
    17 = @("0.000008511", "  +0.06%  ")
    18 = @("14.38",       "  +2.31%  ")
    19 = @("0.9995",      "  -0.05%  ")
    20 = @("27.176.22",   "  +1.09%  ")
    21 = @("5.060",       "  +0.71%  ")
    22 = @("2.149.48",    $null)
    23 = @("10.60",       "  +2.74%  ")
    24 = @("6.441",       "  -0.07%  ")
    25 = @("145.31",      "  -0.76%  ")
    26 = @("1.792",       "  -0.22%  ")
    27 = @("18.09",       "  +1.60%  ")
    28 = @("2.157",       "  +5.03%  ")
    29 = @("114.67",      "  +1.62%  ")
    30 = @("5.000",       "  +7.36%  ")
    31 = @("4.825",       "  +4.06%  ")
    32 = @("0.09224",     "  +0.25%  ")
    33 = @("0.05060",     "  -0.61%  ")
    34 = @("1.197",       "  +4.59%  ")
    35 = @("0.7599",      "  +3.41%  ")
    36 = @("3.027",       "  +1.08%  ")
    37 = @("3.277",       "  +1.44%  ")
    38 = @("2.563",       "  +2.56%  ")
    39 = @("0.5645",      "  +6.23%  ")
    40 = @($null,         "  +0.15%  ")
    41 = @("1.079",       "  +0.62%  ")
    43 = @("118.71",      "  -0.17%  ")
    45 = @("0.1511",      "  +3.07%  ")
    46 = @("0.4824",      "  +4.20%  ")
    47 = @("10.17",       "  +2.76%  ")
    48 = @("0.9999",      "  -0.02%  ")
    49 = @($null,         "  +1.42%  ")
    50 = @("37.19",       "  +0.87%  ")
    51 = @("63.59",       "  +1.21%  ")
}

foreach ($r in $rows.Keys) {
    $pair = $rows[$r]
    $price = $pair[0]
    $volume = $pair[1]
    if ($null -ne $price) {
        Set-TextValue $ws.Cells.Item($r, 4) $price
    }
    if ($null -ne $volume) {
        Set-TextValue $ws.Cells.Item($r, 5) $volume
    }
}

# Row 42 / Row 44: FraxShare and Aptos swap places, each with freshly
# scraped Price/Volume figures.
Set-TextValue $ws.Cells.Item(42, 2) "Aptos"
Set-TextValue $ws.Cells.Item(42, 3) "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Cells.Item(42, 4) "8.928"
Set-TextValue $ws.Cells.Item(42, 5) "  +6.87%  "

Set-TextValue $ws.Cells.Item(44, 2) "FraxShare"
Set-TextValue $ws.Cells.Item(44, 3) "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Cells.Item(44, 4) "6.597"
Set-TextValue $ws.Cells.Item(44, 5) "  +1.83%  "
